# Insert a new data row above row 60 (shifting existing rows 60-150 down to 61-151)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(60).Insert()

$ws.Range("A60").Value = 4
$ws.Range("B60").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C60").Value = "Los Lagos"
$ws.Range("D60").Value = 44897
$ws.Range("E60").Value = 10
$ws.Range("F60").Value = 100112052
$ws.Range("G60").Value = "Albahaca"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 90
$ws.Range("K60").Value = 7000
$ws.Range("L60").Value = 7000
$ws.Range("M60").Value = 7000
$ws.Range("N60").Value = "`$/docena de matas"
$ws.Range("O60").Value = "Región Metropolitana"
$ws.Range("P60").Value = 1167
$ws.Range("Q60").Value = 6
$ws.Range("R60").Value = "Hortaliza"
